$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the 2025 row (row 7) with the refreshed revenue figures
$ws.Range("B7").Value = 2491414.55
$ws.Range("C7").Value = -43.92593416745957
$ws.Range("D7").Value = 2533
$ws.Range("E7").Value = 2533
$ws.Range("F7").Value = 983.5825305961311
$ws.Range("G7").Value = 4.842785543984007
